$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that changed from 45188 (2023-09-19)
# to 45189 (2023-09-20) for every data row (rows 2 through 181).
$ws.Range("C2:C181").Value = 45189
